# Update the DRC country name in the sadc_countries reference sheet:
# "Congo" -> "DRC" (row 5, column C - the Country column)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sadc_countries")

$ws.Range("C5").Value = "DRC"

# Move the active selection to C9, matching the final cursor position
# recorded in the saved workbook.
$ws.Activate()
$ws.Range("C9").Select()
